$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (J) to the table, mirroring the formatting of
# the existing "2020" column (I) for rows 4 (header) through 14 (footer).
$ws.Range("I4:I14").Copy()
$null = $ws.Range("J4:J14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New column values
$ws.Range("J4").Value = 2021
$ws.Range("J5").Value = 1.5
$ws.Range("J6").Value = 0.3
$ws.Range("J7").Value = 0.8
$ws.Range("J8").Value = 0.6
$ws.Range("J9").Value = 1.8
$ws.Range("J10").Value = 0.5
$ws.Range("J11").Value = 0.8
$ws.Range("J12").Value = 1.9
$ws.Range("J13").Value = 4.4000000000000004
$ws.Range("J14").Value = 0.4

# Match the author's final cursor position
$null = $ws.Range("L10").Select()
